# Auto-generated Excel COM-interop script
# Applies the Step1_Data raw-signal update and propagates the
# derived cumulative-sum (Step2_Sj) and threshold-crossing
# (Step3_DataPts_0.5/0.7/0.8/0.9) sheets to match.
$wb = $excel.ActiveWorkbook

# ---- Step1_Data: updated raw per-segment signal values ----
$ws = $wb.Worksheets.Item("Step1_Data")
$step1Changes = @{
    2 = @{ 5 = 0.19844101384188101; 6 = 0.00404520567460941; 7 = 0.21824245874838549; 8 = 0.00813501427056156; 9 = 0.02714665993129684; 10 = 0.00010042143928402; 12 = 0.03971227565587906; 13 = 0.02173438173318116; 14 = 0.04022679127191991; 15 = 0.1069243174595208; 18 = 0.01313637275078756; 20 = 0.00346444475685151; 21 = 0.07405292773869342; 22 = 0.02292536115362107; 23 = 0.00739269029262101; 25 = 0.01573314716786425; 27 = 0.00726796833978142; 28 = 0.01013844062213755; 30 = 0.01877223624850935; 31 = 0.09259629595310405; 32 = 0.03194409060092473; 33 = 0.01181262047321262; 34 = 0.00245803093378806; 36 = 0.02359683294158397 }
    3 = @{ 6 = 0.1818436306629555; 7 = 0.1571337313498879; 8 = 0.02769152360987404; 10 = 0.02055665029256644; 13 = 0.23215479589419671; 14 = 0.0140144927178688; 15 = 0.01177533543614157; 16 = 0.04639017951874923; 17 = 0.00849675432845426; 20 = 0.00110446930619109; 21 = 0.07121175456806542; 22 = 0.02084661091113641; 23 = 0.02052675148585813; 25 = 0.03657247118158554; 28 = 0.02213600738771702; 29 = 0.00413712646633166; 31 = 0.04953936527086527; 32 = 0.04176551174162012; 33 = 0.01272979880011342; 34 = 0.00630305318848765; 35 = 0.00749174881410911; 36 = 0.00218769054849128; 37 = 0.00339054651873347 }
    4 = @{ 5 = 0.08938778970935579; 6 = 0.11086688129066651; 8 = 0.03038863642772606; 9 = 0.04731058117687412; 10 = 0.00577079410181307; 11 = 0.00915988038781556; 12 = 0.1148807999525515; 13 = 0.1884483984738578; 14 = 0.00462112154854733; 15 = 0.00463063940235242; 16 = 0.0032423017389048; 20 = 0.1034507460749499; 21 = 0.08006571263410714; 22 = 0.00726947294129479; 24 = 0.02195798416557627; 25 = 0.002550084726079; 26 = 0.01111472553455613; 27 = 0.0170027290668908; 28 = 0.0144216396253273; 29 = 0.00184261130534627; 30 = 0.0409053270478151; 31 = 0.03607428611688242; 32 = 0.03426155032233501; 33 = 0.00510923820172904; 34 = 0.01267411110530805; 36 = 0.0025919569213376 }
    5 = @{ 4 = 0.02559781919149423; 5 = 0.07059060524737459; 6 = 0.0660677988147827; 7 = 0.02993957147328983; 8 = 0.16822986744847779; 9 = 0.00792534877321889; 10 = 0.00552974469726302; 11 = 0.00140486655319237; 12 = 0.04607746788490362; 13 = 0.17477549059850761; 14 = 0.04436238267960241; 15 = 0.00050032352722895; 17 = 0.01754929960735539; 20 = 0.05248217792455945; 21 = 0.03244302994550256; 22 = 0.03319897575969487; 23 = 0.00275447711611055; 24 = 0.00044783936197227; 27 = 0.01556693450945085; 28 = 0.00481281847096828; 30 = 0.01760823143237504; 31 = 0.04274407467366392; 32 = 0.05883799981861448; 33 = 0.03694703698708807; 34 = 0.0433238660031793; 35 = 0.00028195150012903 }
    6 = @{ 4 = 0.02488905671086378; 5 = 0.01283792850996231; 7 = 0.02798130687488803; 8 = 0.31665766120810401; 10 = 0.01301174448968005; 12 = 0.07719983692295267; 13 = 0.141091528272143; 14 = 0.034080042646624; 15 = 0.01858584719100423; 17 = 0.02237511556020957; 20 = 0.06673517534885963; 21 = 0.03827085704726659; 22 = 0.03409192425055174; 27 = 0.01339719506846298; 28 = 0.0002205999711917; 30 = 0.01890163529232842; 31 = 0.04000249417933899; 32 = 0.04791300019818386; 33 = 0.0214200772889738; 34 = 0.03033697296841052 }
    7 = @{ 5 = 0.19932191828298901; 6 = 0.02913685155973766; 7 = 0.12734413155662411; 8 = 0.00654430093050562; 9 = 0.03080218924774255; 10 = 0.01323606792975902; 11 = 0.00405895694145289; 12 = 0.08183221525228253; 13 = 0.06407547029560899; 14 = 0.019280749713959; 15 = 0.08709257149102653; 17 = 0.00463010731484472; 20 = 0.03162755162479815; 21 = 0.07025305336661679; 22 = 0.02714855477116016; 23 = 0.00377389137609026; 25 = 0.01471573073691566; 27 = 0.02324487209297227; 28 = 0.00458962715140534; 30 = 0.03386020458042406; 31 = 0.0860050351452885; 32 = 0.02332508882738519; 33 = 0.00238168135792826; 36 = 0.01171917845248279 }
    8 = @{ 4 = 0.0027878935787167; 5 = 0.1231138098839325; 6 = 0.126523676913088; 7 = 0; 8 = 0.03202356900343944; 9 = 0.00633505456910358; 12 = 0.28079977803856992; 13 = 0.00123899122631803; 14 = 0.00446514196331527; 15 = 0.12983855468087999; 19 = 0; 20 = 0.03849992958454963; 21 = 0.07864949830954746; 22 = 0; 27 = 0; 29 = 0; 30 = 0.09991355218357642; 31 = 0.07581055006496287; 32 = 0 }
    9 = @{ 4 = 0.00931851467148857; 5 = 0.14901106434673619; 6 = 0.1026638264309416; 7 = 0.0358902875849345; 8 = 0.03971300140275529; 9 = 0.03825481533414943; 10 = 0.00308333210511255; 11 = 0.0080791296129844; 12 = 0.14450784087653529; 13 = 0.03717060597820177; 14 = 0.01216309479668915; 15 = 0.12391035202612551; 16 = 0.00482047504957849; 20 = 0.03668187695858847; 21 = 0.06003395542836706; 22 = 0.02515707275461623; 25 = 0.00299486181941642; 27 = 0.02064096836195304; 29 = 0.00577953816923465; 30 = 0.05914628059222644; 31 = 0.06527706095308845; 32 = 0.01111612443583154; 35 = 0.00458592031044502 }
    10 = @{ 5 = 0.1573657228656474; 6 = 0.1057821027733479; 7 = 0.01850379900145572; 8 = 0.02751520946600475; 9 = 0.04600923688284865; 12 = 0.17222089553732589; 13 = 0.05614301285459129; 14 = 0.00302533057694368; 15 = 0.1555123179385979; 16 = 0; 19 = 0.00309929032481564; 20 = 0.02805691009604994; 21 = 0.05215006988203307; 22 = 0.03728132718516838; 27 = 0.00124024819933785; 30 = 0.05875949257041024; 31 = 0.06492364734926986; 32 = 0.01241138649615182 }
    11 = @{ 4 = 0.18056085956799911; 6 = 0.130610774916193; 7 = 0.01542694759171467; 8 = 0.06292122745800829; 10 = 0.00003610214565773; 11 = 0; 12 = 0.2473917998254172; 14 = 0.11942762232487281; 15 = 0.01317978580464898; 16 = 0; 19 = 0.00188439834982207; 20 = 0.03250302478191167; 21 = 0.06401911929448281; 22 = 0; 27 = 0; 29 = 0.01043023257750334; 30 = 0.08358223941473619; 31 = 0.03802586594703228; 32 = 0 }
}
foreach ($r in $step1Changes.Keys) {
    $rowMap = $step1Changes[$r]
    foreach ($c in $rowMap.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowMap[$c]
    }
}

# ---- Step2_Sj: row-wise cumulative sum of Step1_Data (B:AK) ----
$ws = $wb.Worksheets.Item("Step2_Sj")
$step2Changes = @{
    2 = @{ 5 = 0.19844101384188101; 6 = 0.20248621951649043; 7 = 0.42072867826487592; 8 = 0.42886369253543749; 9 = 0.45601035246673433; 10 = 0.45611077390601834; 11 = 0.45611077390601834; 12 = 0.49582304956189738; 13 = 0.51755743129507858; 14 = 0.55778422256699844; 15 = 0.66470854002651925; 16 = 0.66470854002651925; 17 = 0.66470854002651925; 18 = 0.67784491277730685; 19 = 0.67784491277730685; 20 = 0.68130935753415833; 21 = 0.75536228527285176; 22 = 0.77828764642647286; 23 = 0.78568033671909387; 24 = 0.78568033671909387; 25 = 0.8014134838869581; 26 = 0.8014134838869581; 27 = 0.80868145222673948; 28 = 0.81881989284887702; 29 = 0.81881989284887702; 30 = 0.83759212909738634; 31 = 0.93018842505049038; 32 = 0.96213251565141511; 33 = 0.97394513612462774; 34 = 0.97640316705841579; 35 = 0.97640316705841579 }
    3 = @{ 6 = 0.1818436306629555; 7 = 0.3389773620128434; 8 = 0.36666888562271743; 9 = 0.36666888562271743; 10 = 0.38722553591528386; 11 = 0.38722553591528386; 12 = 0.38722553591528386; 13 = 0.6193803318094806; 14 = 0.63339482452734941; 15 = 0.64517015996349103; 16 = 0.69156033948224027; 17 = 0.70005709381069448; 18 = 0.70005709381069448; 19 = 0.70005709381069448; 20 = 0.7011615631168856; 21 = 0.77237331768495099; 22 = 0.79321992859608736; 23 = 0.81374668008194551; 24 = 0.81374668008194551; 25 = 0.85031915126353108; 26 = 0.85031915126353108; 27 = 0.85031915126353108; 28 = 0.87245515865124812; 29 = 0.87659228511757981; 30 = 0.87659228511757981; 31 = 0.92613165038844503; 32 = 0.96789716213006516; 33 = 0.98062696093017854; 34 = 0.98693001411866621; 35 = 0.99442176293277529; 36 = 0.99660945348126662 }
    4 = @{ 5 = 0.08938778970935579; 6 = 0.20025467100002231; 7 = 0.20025467100002231; 8 = 0.23064330742774836; 9 = 0.27795388860462245; 10 = 0.28372468270643553; 11 = 0.29288456309425109; 12 = 0.40776536304680261; 13 = 0.59621376152066041; 14 = 0.60083488306920774; 15 = 0.60546552247156016; 16 = 0.60870782421046499; 17 = 0.60870782421046499; 18 = 0.60870782421046499; 19 = 0.60870782421046499; 20 = 0.71215857028541485; 21 = 0.79222428291952196; 22 = 0.79949375586081672; 23 = 0.79949375586081672; 24 = 0.82145174002639298; 25 = 0.82400182475247197; 26 = 0.83511655028702814; 27 = 0.85211927935391896; 28 = 0.86654091897924623; 29 = 0.86838353028459248; 30 = 0.90928885733240761; 31 = 0.94536314344929007; 32 = 0.97962469377162509; 33 = 0.98473393197335413; 34 = 0.99740804307866215; 35 = 0.99740804307866215 }
    5 = @{ 4 = 0.02559781919149423; 5 = 0.09618842443886882; 6 = 0.16225622325365152; 7 = 0.19219579472694134; 8 = 0.36042566217541916; 9 = 0.36835101094863804; 10 = 0.37388075564590106; 11 = 0.37528562219909345; 12 = 0.42136309008399708; 13 = 0.59613858068250469; 14 = 0.64050096336210705; 15 = 0.641001286889336; 16 = 0.641001286889336; 17 = 0.6585505864966914; 18 = 0.6585505864966914; 19 = 0.6585505864966914; 20 = 0.71103276442125085; 21 = 0.74347579436675337; 22 = 0.77667477012644826; 23 = 0.7794292472425588; 24 = 0.77987708660453103; 25 = 0.77987708660453103; 26 = 0.77987708660453103; 27 = 0.79544402111398194; 28 = 0.8002568395849502; 29 = 0.8002568395849502; 30 = 0.8178650710173252; 31 = 0.86060914569098912; 32 = 0.9194471455096036; 33 = 0.95639418249669172; 34 = 0.999718048499871 }
    6 = @{ 4 = 0.02488905671086378; 5 = 0.03772698522082609; 6 = 0.03772698522082609; 7 = 0.06570829209571412; 8 = 0.38236595330381812; 9 = 0.38236595330381812; 10 = 0.39537769779349818; 11 = 0.39537769779349818; 12 = 0.47257753471645086; 13 = 0.61366906298859392; 14 = 0.64774910563521793; 15 = 0.66633495282622213; 16 = 0.66633495282622213; 17 = 0.68871006838643167; 18 = 0.68871006838643167; 19 = 0.68871006838643167; 20 = 0.75544524373529132; 21 = 0.79371610078255794; 22 = 0.82780802503310968; 23 = 0.82780802503310968; 24 = 0.82780802503310968; 25 = 0.82780802503310968; 26 = 0.82780802503310968; 27 = 0.84120522010157261; 28 = 0.84142582007276434; 29 = 0.84142582007276434; 30 = 0.86032745536509281; 31 = 0.9003299495444318; 32 = 0.94824294974261569; 33 = 0.96966302703158946 }
    7 = @{ 5 = 0.19932191828298901; 6 = 0.22845876984272667; 7 = 0.35580290139935078; 8 = 0.3623472023298564; 9 = 0.39314939157759893; 10 = 0.40638545950735794; 11 = 0.41044441644881086; 12 = 0.49227663170109337; 13 = 0.55635210199670238; 14 = 0.57563285171066136; 15 = 0.66272542320168792; 16 = 0.66272542320168792; 17 = 0.66735553051653262; 18 = 0.66735553051653262; 19 = 0.66735553051653262; 20 = 0.69898308214133076; 21 = 0.76923613550794756; 22 = 0.79638469027910774; 23 = 0.80015858165519804; 24 = 0.80015858165519804; 25 = 0.81487431239211372; 26 = 0.81487431239211372; 27 = 0.83811918448508593; 28 = 0.84270881163649125; 29 = 0.84270881163649125; 30 = 0.8765690162169153; 31 = 0.96257405136220375; 32 = 0.98589914018958891; 33 = 0.98828082154751717; 34 = 0.98828082154751717; 35 = 0.98828082154751717 }
    8 = @{ 4 = 0.0027878935787167; 5 = 0.1259017034626492; 6 = 0.2524253803757372; 7 = 0.2524253803757372; 8 = 0.28444894937917664; 9 = 0.29078400394828019; 10 = 0.29078400394828019; 11 = 0.29078400394828019; 12 = 0.57158378198685011; 13 = 0.5728227732131681; 14 = 0.57728791517648337; 15 = 0.70712646985736338; 16 = 0.70712646985736338; 17 = 0.70712646985736338; 18 = 0.70712646985736338; 19 = 0.70712646985736338; 20 = 0.745626399441913; 21 = 0.82427589775146048; 22 = 0.82427589775146048; 23 = 0.82427589775146048; 24 = 0.82427589775146048; 25 = 0.82427589775146048; 26 = 0.82427589775146048; 27 = 0.82427589775146048; 28 = 0.82427589775146048; 29 = 0.82427589775146048; 30 = 0.92418944993503693; 31 = 0.99999999999999978 }
    9 = @{ 4 = 0.00931851467148857; 5 = 0.15832957901822475; 6 = 0.26099340544916638; 7 = 0.29688369303410089; 8 = 0.33659669443685619; 9 = 0.37485150977100562; 10 = 0.37793484187611814; 11 = 0.38601397148910255; 12 = 0.53052181236563789; 13 = 0.56769241834383966; 14 = 0.57985551314052886; 15 = 0.7037658651666544; 16 = 0.7085863402162329; 17 = 0.7085863402162329; 18 = 0.7085863402162329; 19 = 0.7085863402162329; 20 = 0.74526821717482139; 21 = 0.80530217260318848; 22 = 0.83045924535780469; 23 = 0.83045924535780469; 24 = 0.83045924535780469; 25 = 0.83345410717722113; 26 = 0.83345410717722113; 27 = 0.85409507553917419; 28 = 0.85409507553917419; 29 = 0.85987461370840879; 30 = 0.91902089430063527; 31 = 0.98429795525372377; 32 = 0.99541407968955531; 33 = 0.99541407968955531; 34 = 0.99541407968955531 }
    10 = @{ 5 = 0.1573657228656474; 6 = 0.2631478256389953; 7 = 0.28165162464045101; 8 = 0.30916683410645579; 9 = 0.35517607098930443; 10 = 0.35517607098930443; 11 = 0.35517607098930443; 12 = 0.52739696652663026; 13 = 0.5835399793812216; 14 = 0.58656530995816525; 15 = 0.74207762789676313; 16 = 0.74207762789676313; 17 = 0.74207762789676313; 18 = 0.74207762789676313; 19 = 0.74517691822157872; 20 = 0.77323382831762866; 21 = 0.82538389819966174; 22 = 0.86266522538483015; 23 = 0.86266522538483015; 24 = 0.86266522538483015; 25 = 0.86266522538483015; 26 = 0.86266522538483015; 27 = 0.86390547358416803; 28 = 0.86390547358416803; 29 = 0.86390547358416803; 30 = 0.92266496615457827; 31 = 0.98758861350384808 }
    11 = @{ 4 = 0.18056085956799911; 5 = 0.18056085956799911; 6 = 0.31117163448419211; 7 = 0.32659858207590675; 8 = 0.38951980953391507; 9 = 0.38951980953391507; 10 = 0.38955591167957282; 11 = 0.38955591167957282; 12 = 0.63694771150498997; 13 = 0.63694771150498997; 14 = 0.75637533382986277; 15 = 0.76955511963451173; 16 = 0.76955511963451173; 17 = 0.76955511963451173; 18 = 0.76955511963451173; 19 = 0.77143951798433386; 20 = 0.80394254276624555; 21 = 0.86796166206072833; 22 = 0.86796166206072833; 23 = 0.86796166206072833; 24 = 0.86796166206072833; 25 = 0.86796166206072833; 26 = 0.86796166206072833; 27 = 0.86796166206072833; 28 = 0.86796166206072833; 29 = 0.87839189463823164; 30 = 0.96197413405296783; 31 = 1 }
}
foreach ($r in $step2Changes.Keys) {
    $rowMap = $step2Changes[$r]
    foreach ($c in $rowMap.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowMap[$c]
    }
}

# ---- Step3_DataPts_0.5: threshold-crossing stats (threshold 0.5) ----
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$step3Changes = @{
    2 = @{ 4 = 12; 6 = 0.51755743129507858; 7 = 10 }
    3 = @{ 6 = 0.6193803318094806 }
    4 = @{ 6 = 0.59621376152066041 }
    5 = @{ 6 = 0.59613858068250469 }
    6 = @{ 5 = 0.03772698522082609; 6 = 0.61366906298859392 }
    7 = @{ 4 = 12; 6 = 0.55635210199670238; 7 = 10 }
    8 = @{ 6 = 0.57158378198685011 }
    9 = @{ 6 = 0.53052181236563789 }
    10 = @{ 6 = 0.52739696652663026 }
    11 = @{ 6 = 0.63694771150498997 }
}
foreach ($r in $step3Changes.Keys) {
    $rowMap = $step3Changes[$r]
    foreach ($c in $rowMap.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowMap[$c]
    }
}

# ---- Step3_DataPts_0.7: threshold-crossing stats (threshold 0.7) ----
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$step4Changes = @{
    2 = @{ 4 = 20; 6 = 0.75536228527285176; 7 = 18 }
    3 = @{ 4 = 16; 6 = 0.70005709381069448; 7 = 13 }
    4 = @{ 6 = 0.71215857028541485 }
    5 = @{ 6 = 0.71103276442125085 }
    6 = @{ 4 = 19; 5 = 0.03772698522082609; 6 = 0.75544524373529132; 7 = 15 }
    7 = @{ 4 = 20; 6 = 0.76923613550794756; 7 = 18 }
    8 = @{ 4 = 14; 6 = 0.70712646985736338; 7 = 12 }
    9 = @{ 6 = 0.7037658651666544 }
    10 = @{ 6 = 0.74207762789676313 }
    11 = @{ 6 = 0.75637533382986277 }
}
foreach ($r in $step4Changes.Keys) {
    $rowMap = $step4Changes[$r]
    foreach ($c in $rowMap.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowMap[$c]
    }
}

# ---- Step3_DataPts_0.8: threshold-crossing stats (threshold 0.8) ----
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$step5Changes = @{
    2 = @{ 4 = 24; 6 = 0.8014134838869581; 7 = 22 }
    3 = @{ 4 = 22; 6 = 0.81374668008194551; 7 = 19 }
    4 = @{ 4 = 23; 6 = 0.82145174002639298; 7 = 21 }
    5 = @{ 4 = 27; 6 = 0.8002568395849502; 7 = 25 }
    6 = @{ 4 = 21; 5 = 0.03772698522082609; 6 = 0.82780802503310968; 7 = 17 }
    7 = @{ 4 = 22; 6 = 0.80015858165519804; 7 = 20 }
    8 = @{ 4 = 20; 6 = 0.82427589775146048; 7 = 18 }
    9 = @{ 6 = 0.80530217260318848 }
    10 = @{ 6 = 0.82538389819966174 }
    11 = @{ 4 = 19; 6 = 0.80394254276624555; 7 = 18 }
}
foreach ($r in $step5Changes.Keys) {
    $rowMap = $step5Changes[$r]
    foreach ($c in $rowMap.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowMap[$c]
    }
}

# ---- Step3_DataPts_0.9: threshold-crossing stats (threshold 0.9) ----
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$step6Changes = @{
    2 = @{ 6 = 0.93018842505049038 }
    3 = @{ 4 = 30; 6 = 0.92613165038844503; 7 = 27 }
    4 = @{ 6 = 0.90928885733240761 }
    5 = @{ 6 = 0.9194471455096036 }
    6 = @{ 5 = 0.03772698522082609; 6 = 0.9003299495444318 }
    7 = @{ 6 = 0.96257405136220375 }
    8 = @{ 6 = 0.92418944993503693 }
    9 = @{ 6 = 0.91902089430063527 }
    10 = @{ 6 = 0.92266496615457827 }
    11 = @{ 6 = 0.96197413405296783 }
}
foreach ($r in $step6Changes.Keys) {
    $rowMap = $step6Changes[$r]
    foreach ($c in $rowMap.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowMap[$c]
    }
}

